# Apply updated crypto price/volume data as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.407.74'
$ws.Range('E2').Value = '  -0.90%  '
$ws.Range('D3').Value = '1.823.66'
$ws.Range('E3').Value = '  -2.02%  '
$ws.Range('D4').Value = '''1.007'
$ws.Range('E4').Value = '  -0.52%  '
$ws.Range('D5').Value = '''331.43'
$ws.Range('E5').Value = '  -0.38%  '
$ws.Range('D6').Value = '''1.006'
$ws.Range('E6').Value = '  -0.36%  '
$ws.Range('D7').Value = '''0.4550'
$ws.Range('E7').Value = '  -2.08%  '
$ws.Range('D8').Value = '''0.3809'
$ws.Range('E8').Value = '  -2.01%  '
$ws.Range('D9').Value = '''46.24'
$ws.Range('E9').Value = '  -0.10%  '
$ws.Range('D10').Value = '''0.07880'
$ws.Range('E10').Value = '  -0.99%  '
$ws.Range('D11').Value = '''0.9675'
$ws.Range('E11').Value = '  -2.72%  '
$ws.Range('D12').Value = '''20.94'
$ws.Range('E12').Value = '  -2.52%  '
$ws.Range('D13').Value = '1.839.78'
$ws.Range('E13').Value = '  -1.61%  '
$ws.Range('D14').Value = '''5.856'
$ws.Range('E14').Value = '  -1.92%  '
$ws.Range('D15').Value = '''7.018'
$ws.Range('E15').Value = '  -1.93%  '
$ws.Range('D16').Value = '''1.008'
$ws.Range('E16').Value = '  -0.59%  '
$ws.Range('D17').Value = '''88.61'
$ws.Range('E17').Value = '  +0.75%  '
$ws.Range('E18').Value = '  -0.94%  '
$ws.Range('D19').Value = '''0.00001028'
$ws.Range('E19').Value = '  -1.38%  '
$ws.Range('D20').Value = '''17.13'
$ws.Range('D21').Value = '''1.005'
$ws.Range('E21').Value = '  -0.38%  '
$ws.Range('D22').Value = '27.387.13'
$ws.Range('E22').Value = '  -0.97%  '
$ws.Range('D23').Value = '''5.311'
$ws.Range('E23').Value = '  -2.36%  '
$ws.Range('D24').Value = '''10.75'
$ws.Range('E24').Value = '  -0.84%  '
$ws.Range('D25').Value = '''2.307'
$ws.Range('E25').Value = '  -0.48%  '
$ws.Range('D26').Value = '2.048.83'
$ws.Range('E26').Value = '  -1.92%  '
$ws.Range('D27').Value = '''156.95'
$ws.Range('E27').Value = '  -0.93%  '
$ws.Range('D28').Value = '''19.35'
$ws.Range('E28').Value = '  -1.62%  '
$ws.Range('D29').Value = '''2.056'
$ws.Range('E29').Value = '  -2.63%  '
$ws.Range('D30').Value = '''5.225'
$ws.Range('E30').Value = '  -2.42%  '
$ws.Range('D31').Value = '''117.91'
$ws.Range('E31').Value = '  -2.68%  '
$ws.Range('D32').Value = '''0.9428'
$ws.Range('E32').Value = '  -2.67%  '
$ws.Range('D33').Value = '''0.09287'
$ws.Range('E33').Value = '  -1.46%  '
$ws.Range('D34').Value = '''3.576'
$ws.Range('E34').Value = '  -1.68%  '
$ws.Range('D35').Value = '''5.223'
$ws.Range('E35').Value = '  -1.13%  '
$ws.Range('D36').Value = '''1.313'
$ws.Range('E36').Value = '  -1.13%  '
$ws.Range('D37').Value = '''0.05917'
$ws.Range('E37').Value = '  -1.52%  '
$ws.Range('D38').Value = '''0.02182'
$ws.Range('E38').Value = '  -1.35%  '
$ws.Range('D39').Value = '''1.157'
$ws.Range('E39').Value = '  -3.21%  '
$ws.Range('D40').Value = '''7.982'
$ws.Range('E40').Value = '  -1.76%  '
$ws.Range('D41').Value = '''0.5731'
$ws.Range('E41').Value = '  -2.66%  '
$ws.Range('D42').Value = '''0.1829'
$ws.Range('E42').Value = '  -2.53%  '
$ws.Range('D43').Value = '''9.989'
$ws.Range('E43').Value = '  -2.13%  '
$ws.Range('D44').Value = '''1.260'
$ws.Range('E44').Value = '  +0.55%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '''12.01'
$ws.Range('E45').Value = '  -0.67%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '''0.5440'
$ws.Range('E46').Value = '  -2.84%  '
$ws.Range('D47').Value = '''1.859'
$ws.Range('E47').Value = '  -2.67%  '
$ws.Range('D48').Value = '''0.06617'
$ws.Range('E48').Value = '  -2.19%  '
$ws.Range('D49').Value = '''109.90'
$ws.Range('E49').Value = '  -1.77%  '
$ws.Range('D50').Value = '''1.036'
$ws.Range('E50').Value = '  -2.33%  '
$ws.Range('D51').Value = '''1.006'
$ws.Range('E51').Value = '  -0.33%  '
